$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A57").Value = "2024-09-26 00:00:00"
$ws.Range("B57").Value = 75050
$ws.Range("C57").Value = 10650.98
$ws.Range("D57").Value = 9425.65
$ws.Range("E57").Value = 7.0121
